$d = $word.ActiveDocument

# Update the date in the title paragraph
$d.Content.Find.Execute("2025-02-11 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-12 Wednesday", 2) | Out-Null

$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "36+15=51"
$t.Cell(1,2).Range.Text = "11+36=47"
$t.Cell(1,3).Range.Text = "40+31=71"
$t.Cell(1,4).Range.Text = "98-72=26"
$t.Cell(1,5).Range.Text = "55+22=77"

$t.Cell(2,1).Range.Text = "30+68=98"
$t.Cell(2,2).Range.Text = "83-7=76"
$t.Cell(2,3).Range.Text = "30-1=29"
$t.Cell(2,4).Range.Text = "23+24=47"
$t.Cell(2,5).Range.Text = "82-6=76"

$t.Cell(3,1).Range.Text = "96-47=49"
$t.Cell(3,2).Range.Text = "82-62=20"
$t.Cell(3,3).Range.Text = "18-3=15"
$t.Cell(3,4).Range.Text = "75+2=77"
$t.Cell(3,5).Range.Text = "74-18=56"

$t.Cell(4,1).Range.Text = "85-38=47"
$t.Cell(4,2).Range.Text = "79+3=82"
$t.Cell(4,3).Range.Text = "32+53=85"
$t.Cell(4,4).Range.Text = "66-53=13"
$t.Cell(4,5).Range.Text = "28+51=79"

$t.Cell(5,1).Range.Text = "75-47=28"
$t.Cell(5,2).Range.Text = "2+11=13"
$t.Cell(5,3).Range.Text = "36+21=57"
$t.Cell(5,4).Range.Text = "21+62=83"
$t.Cell(5,5).Range.Text = "75+19=94"

$t.Cell(6,1).Range.Text = "99-82=17"
$t.Cell(6,2).Range.Text = "36-14=22"
$t.Cell(6,3).Range.Text = "85-61=24"
$t.Cell(6,4).Range.Text = "28+49=77"
$t.Cell(6,5).Range.Text = "72-42=30"

$t.Cell(7,1).Range.Text = "33-22=11"
$t.Cell(7,2).Range.Text = "67+10=77"
$t.Cell(7,3).Range.Text = "61-25=36"
$t.Cell(7,4).Range.Text = "0+2=2"
$t.Cell(7,5).Range.Text = "20+48=68"

$t.Cell(8,1).Range.Text = "78-77=1"
$t.Cell(8,2).Range.Text = "2+30=32"
$t.Cell(8,3).Range.Text = "23+42=65"
$t.Cell(8,4).Range.Text = "71-36=35"
$t.Cell(8,5).Range.Text = "33-28=5"

$t.Cell(9,1).Range.Text = "98-55=43"
$t.Cell(9,2).Range.Text = "3+16=19"
$t.Cell(9,3).Range.Text = "16+36=52"
$t.Cell(9,4).Range.Text = "39-2=37"
$t.Cell(9,5).Range.Text = "11+26=37"

$t.Cell(10,1).Range.Text = "13+70=83"
$t.Cell(10,2).Range.Text = "46+20=66"
$t.Cell(10,3).Range.Text = "67-23=44"
$t.Cell(10,4).Range.Text = "41+25=66"
$t.Cell(10,5).Range.Text = "10+1=11"

$t.Cell(11,1).Range.Text = "40-25=15"
$t.Cell(11,2).Range.Text = "5+16=21"
$t.Cell(11,3).Range.Text = "93-16=77"
$t.Cell(11,4).Range.Text = "43-9=34"
$t.Cell(11,5).Range.Text = "94-74=20"

$t.Cell(12,1).Range.Text = "23+64=87"
$t.Cell(12,2).Range.Text = "36+45=81"
$t.Cell(12,3).Range.Text = "99-73=26"
$t.Cell(12,4).Range.Text = "80-44=36"
$t.Cell(12,5).Range.Text = "80+13=93"

$t.Cell(13,1).Range.Text = "43-26=17"
$t.Cell(13,2).Range.Text = "65+12=77"
$t.Cell(13,3).Range.Text = "38+25=63"
$t.Cell(13,4).Range.Text = "17+81=98"
$t.Cell(13,5).Range.Text = "62-26=36"

$t.Cell(14,1).Range.Text = "25-15=10"
$t.Cell(14,2).Range.Text = "76-53=23"
$t.Cell(14,3).Range.Text = "85-11=74"
$t.Cell(14,4).Range.Text = "4+55=59"
$t.Cell(14,5).Range.Text = "40-15=25"

$t.Cell(15,1).Range.Text = "24-11=13"
$t.Cell(15,2).Range.Text = "30+50=80"
$t.Cell(15,3).Range.Text = "6+58=64"
$t.Cell(15,4).Range.Text = "1+56=57"
$t.Cell(15,5).Range.Text = "49+37=86"

$t.Cell(16,1).Range.Text = "0+27=27"
$t.Cell(16,2).Range.Text = "36+15=51"
$t.Cell(16,3).Range.Text = "7+20=27"
$t.Cell(16,4).Range.Text = "41+39=80"
$t.Cell(16,5).Range.Text = "60+2=62"

$t.Cell(17,1).Range.Text = "42+4=46"
$t.Cell(17,2).Range.Text = "80-61=19"
$t.Cell(17,3).Range.Text = "62-62=0"
$t.Cell(17,4).Range.Text = "0+47=47"
$t.Cell(17,5).Range.Text = "39+43=82"

$t.Cell(18,1).Range.Text = "32+36=68"
$t.Cell(18,2).Range.Text = "99-18=81"
$t.Cell(18,3).Range.Text = "34+65=99"
$t.Cell(18,4).Range.Text = "4+74=78"
$t.Cell(18,5).Range.Text = "72-35=37"

$t.Cell(19,1).Range.Text = "24+43=67"
$t.Cell(19,2).Range.Text = "37+46=83"
$t.Cell(19,3).Range.Text = "69+13=82"
$t.Cell(19,4).Range.Text = "78+9=87"
$t.Cell(19,5).Range.Text = "29-22=7"

$t.Cell(20,1).Range.Text = "28+28=56"
$t.Cell(20,2).Range.Text = "5+20=25"
$t.Cell(20,3).Range.Text = "46-33=13"
$t.Cell(20,4).Range.Text = "20+22=42"
$t.Cell(20,5).Range.Text = "57+29=86"
